# Added logic for printing day and night together.
# Builds 5 shift "cards" on the Card Display sheet: a date header (row 2:3,
# merged) sitting above a shift-label strip (row 5, merged), each block
# boxed with a thick border and filled with a color that distinguishes the
# day shift (yellow) from the night shift (blue).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yellow = 4388090   # RGB(250,244,66)  -> FAF442
$blue   = 12549951  # RGB(63,127,191)  -> 3F7FBF

$dateText = "Friday March 4, 2016"

$blocks = @(
  @{ Anchor = "A";  HeadRange = "A2:I3";   DataRange = "A5:I5";   Shift = "Cherry Line Production 7:00AM - 3:30PM";   Fill = $yellow },
  @{ Anchor = "J";  HeadRange = "J2:R3";   DataRange = "J5:R5";   Shift = "Cherry Line Sorting 7:00AM - 3:30PM";      Fill = $yellow },
  @{ Anchor = "S";  HeadRange = "S2:AA3";  DataRange = "S5:AA5";  Shift = "Operations 7:00AM - 3:30PM";               Fill = $yellow },
  @{ Anchor = "AB"; HeadRange = "AB2:AJ3"; DataRange = "AB5:AJ5"; Shift = "Cherry Line Production 4:00PM - 12:30AM";  Fill = $blue },
  @{ Anchor = "AK"; HeadRange = "AK2:AS3"; DataRange = "AK5:AS5"; Shift = "Cherry Line Sorting 4:00PM - 12:30AM";     Fill = $blue }
)

foreach ($b in $blocks) {
  # ---- date header block (rows 2:3) ----
  $head = $ws.Range($b.HeadRange)
  $head.Interior.Color = $b.Fill
  $head.Borders.Item(7).LineStyle = 1
  $head.Borders.Item(7).Weight = 4
  $head.Borders.Item(8).LineStyle = 1
  $head.Borders.Item(8).Weight = 4
  $head.Borders.Item(9).LineStyle = 1
  $head.Borders.Item(9).Weight = 4
  $head.Borders.Item(10).LineStyle = 1
  $head.Borders.Item(10).Weight = 4

  $headAnchor = $ws.Range($b.Anchor + "2")
  $headAnchor.Value = $dateText
  $headAnchor.Font.Bold = $true
  $headAnchor.Font.Size = 16
  $headAnchor.HorizontalAlignment = -4108
  $headAnchor.VerticalAlignment = -4108

  $head.Merge()

  # ---- shift label block (row 5) ----
  $data = $ws.Range($b.DataRange)
  $data.Interior.Color = $b.Fill
  $data.Borders.Item(7).LineStyle = 1
  $data.Borders.Item(7).Weight = 4
  $data.Borders.Item(8).LineStyle = 1
  $data.Borders.Item(8).Weight = 4
  $data.Borders.Item(9).LineStyle = 1
  $data.Borders.Item(9).Weight = 4
  $data.Borders.Item(10).LineStyle = 1
  $data.Borders.Item(10).Weight = 4

  $dataAnchor = $ws.Range($b.Anchor + "5")
  $dataAnchor.Value = $b.Shift
  $dataAnchor.Font.Bold = $true
  $dataAnchor.Font.Size = 11
  $dataAnchor.HorizontalAlignment = -4108
  $dataAnchor.VerticalAlignment = -4108

  $data.Merge()
}

$ws.Range("AK5").Select()
